# Apply updated price / 1h-volume-change figures to the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Some price strings (e.g. "310.02") look like valid numbers to Excel's
    # auto-detection and would otherwise be stored as a Double. Force the
    # cell to Text, assign the literal string, then restore the default
    # "Normal" style so no stray number-format style is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '42.786.32'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.546.47'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '310.02'
$ws.Range("E5").Value = '  -2.73%  '
Set-TextValue $ws.Range("D6") '99.25'
$ws.Range("E6").Value = '  +1.82%  '
Set-TextValue $ws.Range("D7") '0.570'
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("E8").Value = '  -0.03%  '
Set-TextValue $ws.Range("D9") '0.532'
$ws.Range("E9").Value = '  -0.82%  '
Set-TextValue $ws.Range("D10") '35.97'
$ws.Range("E10").Value = '  -1.20%  '
Set-TextValue $ws.Range("D11") '0.0806'
$ws.Range("E11").Value = '  -1.51%  '
Set-TextValue $ws.Range("D12") '7.41'
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").Value = '2.941.00'
$ws.Range("E14").Value = '  +0.24%  '
Set-TextValue $ws.Range("D15") '15.88'
$ws.Range("E15").Value = '  +4.61%  '
$ws.Range("D16").Value = '2.557.15'
$ws.Range("E16").Value = '  +4.53%  '
Set-TextValue $ws.Range("D17") '0.840'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").Value = '42.796.56'
$ws.Range("E18").Value = '  -0.58%  '
Set-TextValue $ws.Range("D19") '6.74'
$ws.Range("E19").Value = '  -1.67%  '
Set-TextValue $ws.Range("D20") '12.42'
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("E21").Value = '  -1.57%  '
Set-TextValue $ws.Range("D22") '69.42'
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("E27").Value = '  +0.04%  '
Set-TextValue $ws.Range("D28") '2.34'
$ws.Range("E28").Value = '  -3.30%  '
Set-TextValue $ws.Range("D29") '40.09'
$ws.Range("E29").Value = '  -1.48%  '
Set-TextValue $ws.Range("D30") '10.10'
$ws.Range("E30").Value = '  -3.56%  '
Set-TextValue $ws.Range("D31") '159.12'
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("E32").Value = '  -3.18%  '
Set-TextValue $ws.Range("D33") '0.0803'
$ws.Range("E33").Value = '  +1.19%  '
Set-TextValue $ws.Range("D34") '3.30'
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("E35").Value = '  -3.78%  '
$ws.Range("E36").Value = '  -3.60%  '
$ws.Range("E37").Value = '  +4.95%  '
Set-TextValue $ws.Range("D38") '18.43'
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  -0.86%  '
Set-TextValue $ws.Range("D41") '22.43'
$ws.Range("E41").Value = '  +1.21%  '
Set-TextValue $ws.Range("D42") '4.12'
$ws.Range("E42").Value = '  +7.06%  '
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("E44").Value = '  -1.37%  '
Set-TextValue $ws.Range("D45") '3.23'
$ws.Range("E45").Value = '  -2.01%  '
$ws.Range("D46").Value = '1.993.50'
$ws.Range("E46").Value = '  -1.25%  '
Set-TextValue $ws.Range("D47") '9.07'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = '2.787.06'
$ws.Range("E48").Value = '  -0.08%  '
Set-TextValue $ws.Range("D49") '81.26'
$ws.Range("E49").Value = '  -3.85%  '
Set-TextValue $ws.Range("D50") '0.193'
$ws.Range("E50").Value = '  +0.28%  '
Set-TextValue $ws.Range("D51") '73.47'
$ws.Range("E51").Value = '  -4.09%  '
